$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-6) get cyclically permuted for columns D, L, M, N, O, P, S:
#   new row 2 <- old row 3
#   new row 3 <- old row 2
#   new row 4 <- old row 6
#   new row 5 <- old row 4
#   new row 6 <- old row 5
# Capture the "old" values first, then write them into their new rows.

$cols = "D","L","M","N","O","P","S"

$old = @{}
foreach ($r in 2..6) {
    $old[$r] = @{}
    foreach ($c in $cols) {
        $old[$r][$c] = $ws.Range("$c$r").Value2
    }
}

$mapping = @{ 2 = 3; 3 = 2; 4 = 6; 5 = 4; 6 = 5 }

foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $old[$srcRow][$c]
    }
}
